$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format column A as text so that the leading zeros in the Run ID values
# (e.g. "00001048") are preserved instead of being interpreted as numbers.
$ws.Range("A4:A5").NumberFormat = "@"

# Row 4
$ws.Range("A4").Value = "00001048"
$ws.Range("B4").Value = 9492066
$ws.Range("C4").Value = 0.4
$ws.Range("D4").Value = 0.0001
$ws.Range("E4").Value = 0.9
$ws.Range("F4").Value = 0.8
$ws.Range("G4").Value = 300
$ws.Range("H4").Value = 300

# Row 5
$ws.Range("A5").Value = "00001056"
$ws.Range("B5").Value = 20388646
$ws.Range("C5").Value = 0.4
$ws.Range("D5").Value = 0.0001
$ws.Range("E5").Value = 0.9
$ws.Range("F5").Value = 0.8
$ws.Range("G5").Value = 300
$ws.Range("H5").Value = 300
